$d = $word.ActiveDocument

# --- Add a default footer to the (only) section, matching the small
# 9pt / sz=18 paragraph-mark formatting used throughout this template ---
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)

# Stamp the paragraph-mark run formatting (emits <w:sz>/<w:szCs> on the
# paragraph's rPr) before any run exists, then materialize a (still
# empty) trailing run so the part actually gets minted as footer1.xml
# and wired up as the section's default footer reference.
$ftr.Range.Font.Size = 9
$ftr.Range.Font.SizeBi = 9

$ftrEnd = $ftr.Range
$ftrEnd.Collapse(0)
$ftrEnd.InsertAfter("")

# --- Shrink the top margin / grow the bottom margin to match the
# standardized signature block layout ---
$d.PageSetup.TopMargin = 69.4488188976378
$d.PageSetup.BottomMargin = 45.35433070866142
